$wb = $excel.ActiveWorkbook
$bom = $wb.Worksheets.Item("BoM")
$dnf = $wb.Worksheets.Item("DNF")

# --- Update the "Component Groups:" count (7 -> 8) on both the BoM and DNF sheets ---
$bom.Range("F2").Value = 8
$dnf.Range("F2").Value = 8

# --- Update the "Component Count:" text (202 (200 SMD/ 2 THT) -> 204 (202 SMD/ 2 THT)) ---
$bom.Range("F3").Value = "204 (202 SMD/ 2 THT)"
$dnf.Range("F3").Value = "204 (202 SMD/ 2 THT)"

# --- Add the new DNF row 10 for the R2/R3 (3K3, 0805) resistor group ---
$dnf.Range("A10").Value = "'2"
$dnf.Range("B10").Value = "Resistor"
$dnf.Range("C10").Value = "R"
$dnf.Range("D10").Value = "Device"
$dnf.Range("E10").Value = "R2 R3"
$dnf.Range("F10").Value = "3K3"
$dnf.Range("G10").Value = "R_0805_2012Metric"
$dnf.Range("H10").Value = "Resistor_SMD"
$dnf.Range("I10").Value = "'2"
$dnf.Range("J10").Value = "'0"
$dnf.Range("K10").Value = " (DNF)"
$dnf.Range("L10").Value = "~"
$dnf.Range("M10").Value = "/"
$dnf.Range("N10").Value = "pedalboard-display(2)"
$dnf.Range("O10").Value = "'117.7500"
$dnf.Range("P10").Value = "'60.0875"
$dnf.Range("Q10").Value = "'90.0000"
$dnf.Range("R10").Value = "bottom"
$dnf.Range("S10").Value = "SMD"
$dnf.Range("T10").Value = "no"
$dnf.Range("U10").Value = "'2.8500"
$dnf.Range("V10").Value = "'1.4000"

# Match the row formatting used by the equivalent DNF-colored row group already
# present on the BoM sheet (row 12), so the new row's cell styles/fills line up
# with the rest of the DNF coloring scheme instead of creating brand-new styles.
$bom.Range("A12:V12").Copy()
$dnf.Range("A10:V10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
